$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Colors (BGR-packed OLE values expected by Font.Color)
$purple = 10498160   # FF7030A0 - slt row
$brown  = 3368601    # FF996633 - sll row

# --- Step 1: set cell VALUES first, in the exact order needed so new shared
# strings get appended as sll(36), slt(37), 101010(38), 00000(39), 11(40) ---
$ws.Range("A12").Value = "sll"
$ws.Range("A11").Value = "slt"
$ws.Range("I11").Value = "'101010"
$ws.Range("B9").Value  = "'00000"
$ws.Range("E12").Value = "'11"

$ws.Range("B11").Value = "'00000"
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = "'01"
$ws.Range("E11").Value = "'00"
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = "'00"

$ws.Range("B12").Value = "'00000"
$ws.Range("D12").Value = "'01"
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = "'00"
$ws.Range("I12").Value = "'000000"
$ws.Range("J12").Value = "'00"

# --- Step 2: apply font colors. The first no-quote cell and first
# quote-prefixed cell for each color establish the new cellXfs entries in the
# order the target workbook has them (purple no-quote, purple quote, brown
# no-quote, brown quote), then the remaining cells reuse those styles. ---
$ws.Range("A11").Font.Color = $purple
$ws.Range("B11").Font.Color = $purple
$ws.Range("A12").Font.Color = $brown
$ws.Range("B12").Font.Color = $brown

$ws.Range("C11").Font.Color = $purple
$ws.Range("D11").Font.Color = $purple
$ws.Range("E11").Font.Color = $purple
$ws.Range("F11").Font.Color = $purple
$ws.Range("G11").Font.Color = $purple
$ws.Range("H11").Font.Color = $purple
$ws.Range("I11").Font.Color = $purple
$ws.Range("J11").Font.Color = $purple

$ws.Range("C12").Font.Color = $brown
$ws.Range("D12").Font.Color = $brown
$ws.Range("E12").Font.Color = $brown
$ws.Range("F12").Font.Color = $brown
$ws.Range("G12").Font.Color = $brown
$ws.Range("H12").Font.Color = $brown
$ws.Range("I12").Font.Color = $brown
$ws.Range("J12").Font.Color = $brown

# --- Selection moves to G15 ---
[void]$ws.Range("G15").Select()
